$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "63.046.95"
Set-TextValue $ws.Range("D3") "3.050.98"
Set-TextValue $ws.Range("E3") "  -1.06%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.16%  "
Set-TextValue $ws.Range("D5") "582.42"
Set-TextValue $ws.Range("E5") "  -1.40%  "
Set-TextValue $ws.Range("D6") "150.93"
Set-TextValue $ws.Range("E6") "  -2.38%  "
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.13%  "
Set-TextValue $ws.Range("E8") "  -2.19%  "
Set-TextValue $ws.Range("D9") "3.051.35"
Set-TextValue $ws.Range("E9") "  -0.93%  "
Set-TextValue $ws.Range("D10") "0.152"
Set-TextValue $ws.Range("E10") "  -2.82%  "
Set-TextValue $ws.Range("D11") "5.78"
Set-TextValue $ws.Range("E11") "  -0.84%  "
Set-TextValue $ws.Range("E12") "  -1.86%  "
Set-TextValue $ws.Range("E13") "  -3.16%  "
Set-TextValue $ws.Range("D14") "35.84"
Set-TextValue $ws.Range("E14") "  -4.20%  "
Set-TextValue $ws.Range("E15") "  +2.09%  "
Set-TextValue $ws.Range("D16") "3.553.30"
Set-TextValue $ws.Range("E16") "  -1.24%  "
Set-TextValue $ws.Range("E17") "  -0.93%  "
Set-TextValue $ws.Range("D18") "63.018.55"
Set-TextValue $ws.Range("E18") "  -0.83%  "
Set-TextValue $ws.Range("D19") "3.047.68"
Set-TextValue $ws.Range("E19") "  -1.15%  "
Set-TextValue $ws.Range("D20") "477.68"
Set-TextValue $ws.Range("E20") "  +0.32%  "
Set-TextValue $ws.Range("D21") "14.23"
Set-TextValue $ws.Range("E21") "  -2.46%  "
Set-TextValue $ws.Range("E22") "  -1.48%  "
Set-TextValue $ws.Range("D23") "7.50"
Set-TextValue $ws.Range("E23") "  -0.30%  "
Set-TextValue $ws.Range("D24") "2.37"
Set-TextValue $ws.Range("E24") "  -1.80%  "
Set-TextValue $ws.Range("E25") "  +0.42%  "
Set-TextValue $ws.Range("E26") "  -2.32%  "
Set-TextValue $ws.Range("D27") "10.52"
Set-TextValue $ws.Range("E27") "  +4.98%  "
Set-TextValue $ws.Range("B29") "FirstDigitalUSD"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D29") "1.00"
Set-TextValue $ws.Range("E29") "  -0.15%  "
Set-TextValue $ws.Range("B30") "NEARProtocol"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D30") "7.27"
Set-TextValue $ws.Range("E30") "  -1.15%  "
Set-TextValue $ws.Range("B31") "PancakeSwap"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "2.64"
Set-TextValue $ws.Range("E31") "  -1.46%  "
Set-TextValue $ws.Range("B32") "ImmutableX"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D32") "2.19"
Set-TextValue $ws.Range("E32") "  +0.70%  "
Set-TextValue $ws.Range("D33") "27.71"
Set-TextValue $ws.Range("E33") "  +2.05%  "
Set-TextValue $ws.Range("E34") "  -3.50%  "
Set-TextValue $ws.Range("E35") "  +0.91%  "
Set-TextValue $ws.Range("D36") "0.0₃0806"
Set-TextValue $ws.Range("E36") "  -4.66%  "
Set-TextValue $ws.Range("D37") "5.86"
Set-TextValue $ws.Range("E37") "  -3.26%  "
Set-TextValue $ws.Range("E38") "  -1.51%  "
Set-TextValue $ws.Range("D39") "3.06"
Set-TextValue $ws.Range("E39") "  -9.36%  "
Set-TextValue $ws.Range("D40") "50.20"
Set-TextValue $ws.Range("E40") "  -0.97%  "
Set-TextValue $ws.Range("D41") "9.10"
Set-TextValue $ws.Range("E41") "  -1.86%  "
Set-TextValue $ws.Range("D42") "424.83"
Set-TextValue $ws.Range("E42") "  -4.31%  "
Set-TextValue $ws.Range("D43") "0.284"
Set-TextValue $ws.Range("E43") "  +0.28%  "
Set-TextValue $ws.Range("E44") "  +2.33%  "
Set-TextValue $ws.Range("D45") "2.831.84"
Set-TextValue $ws.Range("E45") "  +1.13%  "
Set-TextValue $ws.Range("E46") "  -0.62%  "
Set-TextValue $ws.Range("D47") "37.98"
Set-TextValue $ws.Range("E47") "  -5.14%  "
Set-TextValue $ws.Range("D48") "127.99"
Set-TextValue $ws.Range("E48") "  -2.18%  "
Set-TextValue $ws.Range("D50") "25.02"
Set-TextValue $ws.Range("E50") "  -2.54%  "
Set-TextValue $ws.Range("E51") "  -0.65%  "
